# Restore revision: update the "Integer min" value for the R30 rule row
# (row 10) in the Rules sheet from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
